$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value without Excel coercing
# numeric-looking strings (e.g. "39.10") into numbers.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '42.931.73'
Set-TextValue $ws.Range("E2") '  +0.48%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.287.65'
Set-TextValue $ws.Range("E3") '  +1.56%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.04%  '

# Row 5
Set-TextValue $ws.Range("D5") '252.23'
Set-TextValue $ws.Range("E5") '  -0.43%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.630'
Set-TextValue $ws.Range("E6") '  -0.85%  '

# Row 7
Set-TextValue $ws.Range("D7") '73.42'
Set-TextValue $ws.Range("E7") '  +1.96%  '

# Row 8
Set-TextValue $ws.Range("E8") '  +0.03%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.651'
Set-TextValue $ws.Range("E9") '  +1.29%  '

# Row 10
Set-TextValue $ws.Range("D10") '39.10'
Set-TextValue $ws.Range("E10") '  -4.89%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0977'
Set-TextValue $ws.Range("E11") '  +1.32%  '

# Row 12
Set-TextValue $ws.Range("D12") '59.05'
Set-TextValue $ws.Range("E12") '  -0.89%  '

# Row 13
Set-TextValue $ws.Range("D13") '7.42'
Set-TextValue $ws.Range("E13") '  +0.73%  '

# Row 14
Set-TextValue $ws.Range("E14") '  +0.64%  '

# Row 15
Set-TextValue $ws.Range("D15") '2.629.76'
Set-TextValue $ws.Range("E15") '  +1.62%  '

# Row 16
Set-TextValue $ws.Range("D16") '15.23'
Set-TextValue $ws.Range("E16") '  +2.86%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.869'
Set-TextValue $ws.Range("E17") '  -1.92%  '

# Row 18
Set-TextValue $ws.Range("D18") '2.286.40'
Set-TextValue $ws.Range("E18") '  +1.41%  '

# Row 19
Set-TextValue $ws.Range("D19") '42.836.33'
Set-TextValue $ws.Range("E19") '  +0.24%  '

# Row 20
Set-TextValue $ws.Range("E20") '  +2.80%  '

# Row 21
Set-TextValue $ws.Range("D21") '6.28'
Set-TextValue $ws.Range("E21") '  +0.60%  '

# Row 22
Set-TextValue $ws.Range("D22") '72.52'
Set-TextValue $ws.Range("E22") '  -0.58%  '

# Row 23
Set-TextValue $ws.Range("D23") '236.80'
Set-TextValue $ws.Range("E23") '  +0.47%  '

# Row 24
Set-TextValue $ws.Range("E24") '  +6.00%  '

# Row 25
Set-TextValue $ws.Range("D25") '3.91'
Set-TextValue $ws.Range("E25") '  -1.95%  '

# Row 26
Set-TextValue $ws.Range("D26") '11.56'
Set-TextValue $ws.Range("E26") '  -1.48%  '

# Row 27
Set-TextValue $ws.Range("E27") '  -0.26%  '

# Row 28
Set-TextValue $ws.Range("E28") '  -2.04%  '

# Row 29
Set-TextValue $ws.Range("E29") '  -1.05%  '

# Row 30
Set-TextValue $ws.Range("D30") '2.14'
Set-TextValue $ws.Range("E30") '  -3.30%  '

# Row 31
Set-TextValue $ws.Range("D31") '166.96'
Set-TextValue $ws.Range("E31") '  -0.50%  '

# Row 32
Set-TextValue $ws.Range("D32") '21.04'
Set-TextValue $ws.Range("E32") '  +0.27%  '

# Row 33
Set-TextValue $ws.Range("D33") '6.45'
Set-TextValue $ws.Range("E33") '  +5.22%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -3.43%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0826'
Set-TextValue $ws.Range("E35") '  +5.06%  '

# Row 36
Set-TextValue $ws.Range("D36") '30.87'
Set-TextValue $ws.Range("E36") '  +7.32%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.127'
Set-TextValue $ws.Range("E37") '  +1.58%  '

# Row 38
Set-TextValue $ws.Range("D38") '4.58'
Set-TextValue $ws.Range("E38") '  +10.29%  '

# Row 39
Set-TextValue $ws.Range("D39") '4.78'
Set-TextValue $ws.Range("E39") '  +1.70%  '

# Row 40
Set-TextValue $ws.Range("E40") '  -3.52%  '

# Row 41
Set-TextValue $ws.Range("D41") '14.04'
Set-TextValue $ws.Range("E41") '  +12.90%  '

# Row 42
Set-TextValue $ws.Range("E42") '  +1.90%  '

# Row 43
Set-TextValue $ws.Range("D43") '5.90'
Set-TextValue $ws.Range("E43") '  -2.59%  '

# Row 44
Set-TextValue $ws.Range("E44") '  +7.05%  '

# Row 45
Set-TextValue $ws.Range("D45") '9.16'
Set-TextValue $ws.Range("E45") '  +2.56%  '

# Row 46
Set-TextValue $ws.Range("D46") '61.75'
Set-TextValue $ws.Range("E46") '  -4.46%  '

# Row 47
Set-TextValue $ws.Range("E47") '  -1.61%  '

# Row 48
Set-TextValue $ws.Range("E48") '  +1.61%  '

# Row 49
Set-TextValue $ws.Range("D49") '102.00'
Set-TextValue $ws.Range("E49") '  +7.78%  '

# Row 50
Set-TextValue $ws.Range("E50") '  +0.14%  '

# Row 51
Set-TextValue $ws.Range("D51") '1.16'
Set-TextValue $ws.Range("E51") '  -2.70%  '
